$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 14 (shifts the old "total" row 14 -> 15 and the
# old footer row 15 -> 16), the same way Excel's UI "Insert Row" does.
$ws.Rows.Item(14).Insert()

# Bring over the per-column cell formatting from the row-13 item template
# (same style pattern used by every item row 7-13) so row 14 looks like a
# normal item row instead of a blank one.
$ws.Range("A13:B13").Copy()
$ws.Range("A14:B14").PasteSpecial(-4122)

$ws.Range("C13:G13").Copy()
$ws.Range("C14:G14").PasteSpecial(-4122)

$ws.Range("H13:K13").Copy()
$ws.Range("H14:K14").PasteSpecial(-4122)

$ws.Range("L13:M13").Copy()
$ws.Range("L14:M14").PasteSpecial(-4122)

$ws.Range("N13:O13").Copy()
$ws.Range("N14:O14").PasteSpecial(-4122)

$ws.Range("P13").Copy()
$ws.Range("P14").PasteSpecial(-4122)

$ws.Range("Q13").Copy()
$ws.Range("Q14").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# Re-create the merges for the new item row (mirrors rows 7-13).
$ws.Range("A14:B14").Merge()
$ws.Range("C14:G14").Merge()
$ws.Range("H14:K14").Merge()
$ws.Range("L14:M14").Merge()
$ws.Range("N14:O14").Merge()

$ws.Rows.Item(14).RowHeight = 25.5

# New item #8 data.
$ws.Range("A14").Value = 8
$ws.Range("C14").Value = "زيت برافين"
$ws.Range("H14").Value = "2:0"
$ws.Range("L14").Value = "0"
$ws.Range("N14").Value = "20.00"
$ws.Range("P14").Value = "20.0000"
$ws.Range("Q14").Value = "1:0"

# Update the totals row (shifted down to row 15): add the new item's sell
# price (20.00) to the previous total (266.045 -> 286.045).
$ws.Range("P15").Value = 286.04500000000002
$ws.Rows.Item(15).RowHeight = 24.75

# Update the footer timestamp (shifted down to row 16) to the new export time.
$ws.Range("A16").Value = "Monday, 25 August, 2025 11:12 AM"
